# Map032 scene update: move the duplicated "English" column from C to D
# (shifts the C1:C13 data - present only on rows 1,2,3,4,6 - into column D,
# leaving column C empty and extending the sheet's used range to A1:D13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Cut($ws.Range("D1"))
$ws.Range("C2").Cut($ws.Range("D2"))
$ws.Range("C3").Cut($ws.Range("D3"))
$ws.Range("C4").Cut($ws.Range("D4"))
$ws.Range("C6").Cut($ws.Range("D6"))
